$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 10114
$ws1.Range("F8").Value = 1611
$ws1.Range("F22").Value = 317
$ws1.Range("F31").Value = 400
$ws1.Range("F33").Value = 371
$ws1.Range("F35").Value = 603

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F18").Value = 1084
$ws2.Range("F20").Value = 605

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F6").Value = 2514
$ws3.Range("F7").Value = 4061
$ws3.Range("F8").Value = 58

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 10114
$ws4.Range("F7").Value = 4061
$ws4.Range("F8").Value = 58
$ws4.Range("F12").Value = 1611
$ws4.Range("F25").Value = 1084
$ws4.Range("F27").Value = 317
$ws4.Range("F35").Value = 400
$ws4.Range("F37").Value = 371
$ws4.Range("F39").Value = 603
